# Auto-generated edit script: restores the species-record permutation
# described by the diff (each data row receives another row's species-block).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($row, $col, $type, $value) {
    $cell = $ws.Cells.Item($row, $col)
    if ($type -eq "num") {
        $cell.Value = $value
    } elseif ($type -eq "str") {
        $cell.Value = $value
    } elseif ($type -eq "textnum") {
        $cell.NumberFormat = "@"
        $cell.Value = $value
    } else {
        $cell.Value = ""
    }
}

# Row 2
Set-Cell 2 1 "num" 111815501
Set-Cell 2 2 "num" 77515
Set-Cell 2 4 "str" "NT"
Set-Cell 2 5 "num" 6425
Set-Cell 2 6 "str" "Garnlav"
Set-Cell 2 7 "str" "Alectoria sarmentosa"
Set-Cell 2 8 "str" "(Ach.) Ach."
Set-Cell 2 9 "empty" $null
Set-Cell 2 11 "absent" $null
Set-Cell 2 12 "absent" $null
Set-Cell 2 13 "absent" $null
Set-Cell 2 14 "absent" $null
Set-Cell 2 17 "num" 457501.6028477412
Set-Cell 2 18 "num" 7058396.615926245
Set-Cell 2 29 "absent" $null

# Row 3
Set-Cell 3 1 "num" 111815504
Set-Cell 3 2 "num" 77515
Set-Cell 3 4 "str" "NT"
Set-Cell 3 5 "num" 6425
Set-Cell 3 6 "str" "Garnlav"
Set-Cell 3 7 "str" "Alectoria sarmentosa"
Set-Cell 3 8 "str" "(Ach.) Ach."
Set-Cell 3 9 "empty" $null
Set-Cell 3 11 "absent" $null
Set-Cell 3 12 "absent" $null
Set-Cell 3 13 "absent" $null
Set-Cell 3 14 "absent" $null
Set-Cell 3 17 "num" 457647.6599703166
Set-Cell 3 18 "num" 7057607.672130827
Set-Cell 3 29 "absent" $null

# Row 4
Set-Cell 4 1 "num" 111815494
Set-Cell 4 2 "num" 89419
Set-Cell 4 4 "str" "NT"
Set-Cell 4 5 "num" 1204
Set-Cell 4 6 "str" "Gränsticka"
Set-Cell 4 7 "str" "Phellopilus nigrolimitatus"
Set-Cell 4 8 "str" "(Romell) Niemelä, T.Wagner & M.Fisch."
Set-Cell 4 9 "empty" $null
Set-Cell 4 11 "absent" $null
Set-Cell 4 12 "absent" $null
Set-Cell 4 13 "absent" $null
Set-Cell 4 14 "absent" $null
Set-Cell 4 17 "num" 457558.4152710024
Set-Cell 4 18 "num" 7057456.868363639
Set-Cell 4 29 "absent" $null

# Row 5
Set-Cell 5 1 "num" 111815484
Set-Cell 5 2 "num" 56398
Set-Cell 5 4 "str" "NT"
Set-Cell 5 5 "num" 100109
Set-Cell 5 6 "str" "Tretåig hackspett"
Set-Cell 5 7 "str" "Picoides tridactylus"
Set-Cell 5 8 "str" "(Linnaeus, 1758)"
Set-Cell 5 9 "empty" $null
Set-Cell 5 11 "empty" $null
Set-Cell 5 12 "empty" $null
Set-Cell 5 13 "empty" $null
Set-Cell 5 14 "empty" $null
Set-Cell 5 17 "num" 457499.2570064011
Set-Cell 5 18 "num" 7058354.436788658
Set-Cell 5 29 "str" "ringhack"

# Row 6
Set-Cell 6 1 "num" 111815486
Set-Cell 6 2 "num" 56398
Set-Cell 6 4 "str" "NT"
Set-Cell 6 5 "num" 100109
Set-Cell 6 6 "str" "Tretåig hackspett"
Set-Cell 6 7 "str" "Picoides tridactylus"
Set-Cell 6 8 "str" "(Linnaeus, 1758)"
Set-Cell 6 9 "empty" $null
Set-Cell 6 11 "empty" $null
Set-Cell 6 12 "empty" $null
Set-Cell 6 13 "empty" $null
Set-Cell 6 14 "empty" $null
Set-Cell 6 17 "num" 457490.629824138
Set-Cell 6 18 "num" 7057910.64054891
Set-Cell 6 29 "str" "ringhack"

# Row 7
Set-Cell 7 1 "num" 111815478
Set-Cell 7 2 "num" 90087
Set-Cell 7 4 "str" "LC"
Set-Cell 7 5 "num" 3298
Set-Cell 7 6 "str" "Trådticka"
Set-Cell 7 7 "str" "Climacocystis borealis"
Set-Cell 7 8 "str" "(Fr.) Kotl. & Pouzar"
Set-Cell 7 9 "empty" $null
Set-Cell 7 11 "absent" $null
Set-Cell 7 12 "absent" $null
Set-Cell 7 13 "absent" $null
Set-Cell 7 14 "absent" $null
Set-Cell 7 17 "num" 457490.7608241383
Set-Cell 7 18 "num" 7057588.885967719
Set-Cell 7 29 "absent" $null

# Row 8
Set-Cell 8 1 "num" 111815476
Set-Cell 8 2 "num" 90087
Set-Cell 8 4 "str" "LC"
Set-Cell 8 5 "num" 3298
Set-Cell 8 6 "str" "Trådticka"
Set-Cell 8 7 "str" "Climacocystis borealis"
Set-Cell 8 8 "str" "(Fr.) Kotl. & Pouzar"
Set-Cell 8 9 "empty" $null
Set-Cell 8 11 "absent" $null
Set-Cell 8 12 "absent" $null
Set-Cell 8 13 "absent" $null
Set-Cell 8 14 "absent" $null
Set-Cell 8 17 "num" 457561.367873844
Set-Cell 8 18 "num" 7058241.631711838
Set-Cell 8 29 "absent" $null

# Row 9
Set-Cell 9 1 "num" 111815505
Set-Cell 9 2 "num" 77515
Set-Cell 9 4 "str" "NT"
Set-Cell 9 5 "num" 6425
Set-Cell 9 6 "str" "Garnlav"
Set-Cell 9 7 "str" "Alectoria sarmentosa"
Set-Cell 9 8 "str" "(Ach.) Ach."
Set-Cell 9 9 "empty" $null
Set-Cell 9 11 "absent" $null
Set-Cell 9 12 "absent" $null
Set-Cell 9 13 "absent" $null
Set-Cell 9 14 "absent" $null
Set-Cell 9 17 "num" 457627.5754243882
Set-Cell 9 18 "num" 7057503.498196352
Set-Cell 9 29 "absent" $null

# Row 10
Set-Cell 10 1 "num" 111815475
Set-Cell 10 2 "num" 90087
Set-Cell 10 4 "str" "LC"
Set-Cell 10 5 "num" 3298
Set-Cell 10 6 "str" "Trådticka"
Set-Cell 10 7 "str" "Climacocystis borealis"
Set-Cell 10 8 "str" "(Fr.) Kotl. & Pouzar"
Set-Cell 10 9 "empty" $null
Set-Cell 10 11 "absent" $null
Set-Cell 10 12 "absent" $null
Set-Cell 10 13 "absent" $null
Set-Cell 10 14 "absent" $null
Set-Cell 10 17 "num" 457550.3597693135
Set-Cell 10 18 "num" 7058250.221744461
Set-Cell 10 29 "absent" $null

# Row 11
Set-Cell 11 1 "num" 111815503
Set-Cell 11 2 "num" 77515
Set-Cell 11 4 "str" "NT"
Set-Cell 11 5 "num" 6425
Set-Cell 11 6 "str" "Garnlav"
Set-Cell 11 7 "str" "Alectoria sarmentosa"
Set-Cell 11 8 "str" "(Ach.) Ach."
Set-Cell 11 9 "empty" $null
Set-Cell 11 11 "absent" $null
Set-Cell 11 12 "absent" $null
Set-Cell 11 13 "absent" $null
Set-Cell 11 14 "absent" $null
Set-Cell 11 17 "num" 457482.2929676044
Set-Cell 11 18 "num" 7057720.548935141
Set-Cell 11 29 "absent" $null

# Row 12
Set-Cell 12 1 "num" 111815490
Set-Cell 12 2 "num" 56414
Set-Cell 12 4 "str" "NT"
Set-Cell 12 5 "num" 100049
Set-Cell 12 6 "str" "Spillkråka"
Set-Cell 12 7 "str" "Dryocopus martius"
Set-Cell 12 8 "str" "(Linnaeus, 1758)"
Set-Cell 12 9 "empty" $null
Set-Cell 12 11 "empty" $null
Set-Cell 12 12 "empty" $null
Set-Cell 12 13 "empty" $null
Set-Cell 12 14 "empty" $null
Set-Cell 12 17 "num" 457486.844484477
Set-Cell 12 18 "num" 7058059.55768314
Set-Cell 12 29 "str" "hack"

# Row 13
Set-Cell 13 1 "num" 111815482
Set-Cell 13 2 "num" 56398
Set-Cell 13 4 "str" "NT"
Set-Cell 13 5 "num" 100109
Set-Cell 13 6 "str" "Tretåig hackspett"
Set-Cell 13 7 "str" "Picoides tridactylus"
Set-Cell 13 8 "str" "(Linnaeus, 1758)"
Set-Cell 13 9 "empty" $null
Set-Cell 13 11 "empty" $null
Set-Cell 13 12 "empty" $null
Set-Cell 13 13 "empty" $null
Set-Cell 13 14 "empty" $null
Set-Cell 13 17 "num" 457734.0897740572
Set-Cell 13 18 "num" 7057881.607121572
Set-Cell 13 29 "str" "ringhack"

# Row 14
Set-Cell 14 1 "num" 111815480
Set-Cell 14 2 "num" 90087
Set-Cell 14 4 "str" "LC"
Set-Cell 14 5 "num" 3298
Set-Cell 14 6 "str" "Trådticka"
Set-Cell 14 7 "str" "Climacocystis borealis"
Set-Cell 14 8 "str" "(Fr.) Kotl. & Pouzar"
Set-Cell 14 9 "empty" $null
Set-Cell 14 11 "absent" $null
Set-Cell 14 12 "absent" $null
Set-Cell 14 13 "absent" $null
Set-Cell 14 14 "absent" $null
Set-Cell 14 17 "num" 457650.8748659134
Set-Cell 14 18 "num" 7057581.852142417
Set-Cell 14 29 "absent" $null

# Row 15
Set-Cell 15 1 "num" 111815483
Set-Cell 15 2 "num" 56398
Set-Cell 15 4 "str" "NT"
Set-Cell 15 5 "num" 100109
Set-Cell 15 6 "str" "Tretåig hackspett"
Set-Cell 15 7 "str" "Picoides tridactylus"
Set-Cell 15 8 "str" "(Linnaeus, 1758)"
Set-Cell 15 9 "textnum" "1"
Set-Cell 15 11 "empty" $null
Set-Cell 15 12 "empty" $null
Set-Cell 15 13 "str" "födosökande"
Set-Cell 15 14 "str" "observerad"
Set-Cell 15 17 "num" 457815.3956129756
Set-Cell 15 18 "num" 7058239.570048946
Set-Cell 15 29 "absent" $null

# Row 16
Set-Cell 16 1 "num" 111815495
Set-Cell 16 2 "num" 89423
Set-Cell 16 4 "str" "NT"
Set-Cell 16 5 "num" 5432
Set-Cell 16 6 "str" "Granticka"
Set-Cell 16 7 "str" "Porodaedalea chrysoloma"
Set-Cell 16 8 "str" "(Fr.) Fiasson & Niemelä"
Set-Cell 16 9 "empty" $null
Set-Cell 16 11 "absent" $null
Set-Cell 16 12 "absent" $null
Set-Cell 16 13 "absent" $null
Set-Cell 16 14 "absent" $null
Set-Cell 16 17 "num" 457740.1420321366
Set-Cell 16 18 "num" 7057634.880048735
Set-Cell 16 29 "absent" $null

# Row 17
Set-Cell 17 1 "num" 111815469
Set-Cell 17 2 "num" 90087
Set-Cell 17 4 "str" "LC"
Set-Cell 17 5 "num" 3298
Set-Cell 17 6 "str" "Trådticka"
Set-Cell 17 7 "str" "Climacocystis borealis"
Set-Cell 17 8 "str" "(Fr.) Kotl. & Pouzar"
Set-Cell 17 9 "empty" $null
Set-Cell 17 11 "absent" $null
Set-Cell 17 12 "absent" $null
Set-Cell 17 13 "absent" $null
Set-Cell 17 14 "absent" $null
Set-Cell 17 17 "num" 457736.9978307564
Set-Cell 17 18 "num" 7057632.69988044
Set-Cell 17 29 "absent" $null

# Row 18
Set-Cell 18 1 "num" 111815498
Set-Cell 18 2 "num" 89423
Set-Cell 18 4 "str" "NT"
Set-Cell 18 5 "num" 5432
Set-Cell 18 6 "str" "Granticka"
Set-Cell 18 7 "str" "Porodaedalea chrysoloma"
Set-Cell 18 8 "str" "(Fr.) Fiasson & Niemelä"
Set-Cell 18 9 "empty" $null
Set-Cell 18 11 "absent" $null
Set-Cell 18 12 "absent" $null
Set-Cell 18 13 "absent" $null
Set-Cell 18 14 "absent" $null
Set-Cell 18 17 "num" 457525.8934188869
Set-Cell 18 18 "num" 7057587.081301005
Set-Cell 18 29 "absent" $null

# Row 19
Set-Cell 19 1 "num" 111815500
Set-Cell 19 2 "num" 77515
Set-Cell 19 4 "str" "NT"
Set-Cell 19 5 "num" 6425
Set-Cell 19 6 "str" "Garnlav"
Set-Cell 19 7 "str" "Alectoria sarmentosa"
Set-Cell 19 8 "str" "(Ach.) Ach."
Set-Cell 19 9 "empty" $null
Set-Cell 19 11 "absent" $null
Set-Cell 19 12 "absent" $null
Set-Cell 19 13 "absent" $null
Set-Cell 19 14 "absent" $null
Set-Cell 19 17 "num" 457856.1313392611
Set-Cell 19 18 "num" 7058258.134138036
Set-Cell 19 29 "absent" $null

# Row 20
Set-Cell 20 1 "num" 111815472
Set-Cell 20 2 "num" 90087
Set-Cell 20 4 "str" "LC"
Set-Cell 20 5 "num" 3298
Set-Cell 20 6 "str" "Trådticka"
Set-Cell 20 7 "str" "Climacocystis borealis"
Set-Cell 20 8 "str" "(Fr.) Kotl. & Pouzar"
Set-Cell 20 9 "empty" $null
Set-Cell 20 11 "absent" $null
Set-Cell 20 12 "absent" $null
Set-Cell 20 13 "absent" $null
Set-Cell 20 14 "absent" $null
Set-Cell 20 17 "num" 457859.1684109565
Set-Cell 20 18 "num" 7058252.317324108
Set-Cell 20 29 "absent" $null

# Row 21
Set-Cell 21 1 "num" 111815499
Set-Cell 21 2 "num" 89423
Set-Cell 21 4 "str" "NT"
Set-Cell 21 5 "num" 5432
Set-Cell 21 6 "str" "Granticka"
Set-Cell 21 7 "str" "Porodaedalea chrysoloma"
Set-Cell 21 8 "str" "(Fr.) Fiasson & Niemelä"
Set-Cell 21 9 "empty" $null
Set-Cell 21 11 "absent" $null
Set-Cell 21 12 "absent" $null
Set-Cell 21 13 "absent" $null
Set-Cell 21 14 "absent" $null
Set-Cell 21 17 "num" 457639.6605191349
Set-Cell 21 18 "num" 7057508.669857187
Set-Cell 21 29 "absent" $null

# Row 22
Set-Cell 22 1 "num" 111815506
Set-Cell 22 2 "num" 77515
Set-Cell 22 4 "str" "NT"
Set-Cell 22 5 "num" 6425
Set-Cell 22 6 "str" "Garnlav"
Set-Cell 22 7 "str" "Alectoria sarmentosa"
Set-Cell 22 8 "str" "(Ach.) Ach."
Set-Cell 22 9 "empty" $null
Set-Cell 22 11 "absent" $null
Set-Cell 22 12 "absent" $null
Set-Cell 22 13 "absent" $null
Set-Cell 22 14 "absent" $null
Set-Cell 22 17 "num" 457558.4688635201
Set-Cell 22 18 "num" 7057460.867698954
Set-Cell 22 29 "absent" $null

# Row 23
Set-Cell 23 1 "num" 111815470
Set-Cell 23 2 "num" 90087
Set-Cell 23 4 "str" "LC"
Set-Cell 23 5 "num" 3298
Set-Cell 23 6 "str" "Trådticka"
Set-Cell 23 7 "str" "Climacocystis borealis"
Set-Cell 23 8 "str" "(Fr.) Kotl. & Pouzar"
Set-Cell 23 9 "empty" $null
Set-Cell 23 11 "absent" $null
Set-Cell 23 12 "absent" $null
Set-Cell 23 13 "absent" $null
Set-Cell 23 14 "absent" $null
Set-Cell 23 17 "num" 457615.01761246
Set-Cell 23 18 "num" 7058260.908339346
Set-Cell 23 29 "absent" $null

# Row 24
Set-Cell 24 1 "num" 111815489
Set-Cell 24 2 "num" 56414
Set-Cell 24 4 "str" "NT"
Set-Cell 24 5 "num" 100049
Set-Cell 24 6 "str" "Spillkråka"
Set-Cell 24 7 "str" "Dryocopus martius"
Set-Cell 24 8 "str" "(Linnaeus, 1758)"
Set-Cell 24 9 "empty" $null
Set-Cell 24 11 "empty" $null
Set-Cell 24 12 "empty" $null
Set-Cell 24 13 "empty" $null
Set-Cell 24 14 "empty" $null
Set-Cell 24 17 "num" 457851.1019836199
Set-Cell 24 18 "num" 7058247.981310523
Set-Cell 24 29 "str" "hack"

# Row 25
Set-Cell 25 1 "num" 111815492
Set-Cell 25 2 "num" 78578
Set-Cell 25 4 "str" "NT"
Set-Cell 25 5 "num" 6458
Set-Cell 25 6 "str" "Lunglav"
Set-Cell 25 7 "str" "Lobaria pulmonaria"
Set-Cell 25 8 "str" "(L.) Hoffm."
Set-Cell 25 9 "empty" $null
Set-Cell 25 11 "absent" $null
Set-Cell 25 12 "absent" $null
Set-Cell 25 13 "absent" $null
Set-Cell 25 14 "absent" $null
Set-Cell 25 17 "num" 457652.3413775756
Set-Cell 25 18 "num" 7058423.476722932
Set-Cell 25 29 "absent" $null

# Row 26
Set-Cell 26 1 "num" 111815471
Set-Cell 26 2 "num" 90087
Set-Cell 26 4 "str" "LC"
Set-Cell 26 5 "num" 3298
Set-Cell 26 6 "str" "Trådticka"
Set-Cell 26 7 "str" "Climacocystis borealis"
Set-Cell 26 8 "str" "(Fr.) Kotl. & Pouzar"
Set-Cell 26 9 "empty" $null
Set-Cell 26 11 "absent" $null
Set-Cell 26 12 "absent" $null
Set-Cell 26 13 "absent" $null
Set-Cell 26 14 "absent" $null
Set-Cell 26 17 "num" 457690.9089585465
Set-Cell 26 18 "num" 7058279.888449568
Set-Cell 26 29 "absent" $null

# Row 27
Set-Cell 27 1 "num" 111815485
Set-Cell 27 2 "num" 56398
Set-Cell 27 4 "str" "NT"
Set-Cell 27 5 "num" 100109
Set-Cell 27 6 "str" "Tretåig hackspett"
Set-Cell 27 7 "str" "Picoides tridactylus"
Set-Cell 27 8 "str" "(Linnaeus, 1758)"
Set-Cell 27 9 "empty" $null
Set-Cell 27 11 "empty" $null
Set-Cell 27 12 "empty" $null
Set-Cell 27 13 "empty" $null
Set-Cell 27 14 "empty" $null
Set-Cell 27 17 "num" 457446.9368417656
Set-Cell 27 18 "num" 7058136.079544679
Set-Cell 27 29 "str" "ringhack"

